# Clean up code and fix output
# Adds a new "Yearly demand" worksheet (same layout style as the other
# dispatch sheets: header row 0..23 across B1:Y1, row index 0..2 down
# A2:A4) as the last sheet in the workbook.

$wb = $excel.ActiveWorkbook

# Reference sheet whose header/index formatting (bold, centered, top
# aligned, thin-box border -- style index 1 in the original workbook) we
# want to replicate exactly on the new sheet.
$formatSource = $wb.Worksheets.Item(1)

# New sheet goes after the current last sheet ("Connected Households").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Yearly demand"

# Copy the header/index cell formatting from the reference sheet so the
# new sheet reuses the same style (bold font, centered/top aligned, thin
# border) rather than creating a brand-new style entry. Note: A1 itself
# is never populated/formatted on the reference sheet, so copy the
# header row and index column separately instead of the full A1:Y4 block.
$formatSource.Range("B1:Y1").Copy()
$newSheet.Range("B1:Y1").PasteSpecial(-4122) # xlPasteFormats

$formatSource.Range("A2:A4").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122) # xlPasteFormats

# Header row: 0..23 across B1:Y1
for ($i = 0; $i -le 23; $i++) {
    $col = $i + 2 # B = column 2
    $newSheet.Cells.Item(1, $col).Value = $i
}

# Row index column: 0,1,2 down A2:A4
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2

# Data rows (B..Y => columns 2..25)
$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $newSheet.Cells.Item(2, $i + 2).Value = $row2[$i]
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    $newSheet.Cells.Item(3, $i + 2).Value = $row3[$i]
}
for ($i = 0; $i -lt $row4.Length; $i++) {
    $newSheet.Cells.Item(4, $i + 2).Value = $row4[$i]
}

$newSheet.Range("A1").Select()
